$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "kdj"
$ws.Range("C1").Value = "dk"
$ws.Range("D1").Value = "sk"

$ws.Range("D1").Select() | Out-Null
